# Update the prediction-score column (B) with the newly computed values,
# replacing the placeholder "1"s left over from the previous copy of
# ful-path.csv.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.3185113181394
$ws.Range("B3").Value = 19.418463845647892
$ws.Range("B4").Value = 9.4790332671935289
$ws.Range("B5").Value = 2.8848686594246864
$ws.Range("B6").Value = 7.2698368973468632
$ws.Range("B7").Value = 18.42405118119634
$ws.Range("B8").Value = 20.648375805302543
$ws.Range("B9").Value = 27.781549545810897

# Re-affirm the text formatting on the header row and the row-label column
# (kept as "@" / text format, matching the workbook's existing style).
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A9").NumberFormat = "@"
